$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.473.34"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "1.580.52"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.246"
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = "  -1.53%  "

$ws.Range("E12").Value = "  +0.91%  "

$ws.Range("D13").Value = "1.804.94"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").Value = "1.585.51"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.519"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "28.482.83"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "0.0₃0691"
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("E23").Value = "  -3.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.83%  "

$ws.Range("E25").Value = "  +4.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.76%  "

$ws.Range("E29").Value = "  -2.28%  "

$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0484"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.45%  "

$ws.Range("E32").Value = "  -1.78%  "

$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "

$ws.Range("D35").Value = "1.400.81"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.50%  "

$ws.Range("E37").Value = "  -3.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("E40").Value = "  -0.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.522"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.63%  "

$ws.Range("E42").Value = "  +0.29%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.21%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.793"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0463"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.961"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").Value = "1.718.04"
$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0517"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
